# qs-addandsm.docx: collapse the word-by-word runs in the Title, Author
# and Abstract paragraphs into a single run each (text content is
# unchanged - only the run segmentation changes).

$d = $word.ActiveDocument

# Paragraph 1 (style "Title"): "Questions: Vector addition and scalar multiplication"
$titleText = "Questions: Vector addition and scalar multiplication"
$r = $d.Paragraphs(1).Range
$r.Find.Execute($titleText, $true, $false, $false, $false, $false, $true, 1, $false, $titleText, 2)

# Paragraph 2 (style "Author"): "Renee Knapp, Kin Wang Pang"
$authorText = "Renee Knapp, Kin Wang Pang"
$r = $d.Paragraphs(2).Range
$r.Find.Execute($authorText, $true, $false, $false, $false, $false, $true, 1, $false, $authorText, 2)

# Paragraph 4 (style "Abstract"): "A selection of questions for the study guide on vector addition and scalar multiplication."
$abstractText = "A selection of questions for the study guide on vector addition and scalar multiplication."
$r = $d.Paragraphs(4).Range
$r.Find.Execute($abstractText, $true, $false, $false, $false, $false, $true, 1, $false, $abstractText, 2)
